# Update the "test mode template" zip code sample values on the Zipcode sheet
# and move the active selection to the cell that was edited (B4).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zipcode")

$ws.Range("B4").Value = 43215
$ws.Range("B6").Value = 94043

$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
